$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2,8).Value = 379.75  # ALC!H2 was 645.1818
$ws.Cells.Item(2,10).Value = 399.66666  # ALC!J2 was 876.3333
$ws.Cells.Item(2,12).Value = 399.66666  # ALC!L2 was 876.3333
$ws.Cells.Item(2,14).Value = -625.66666  # ALC!N2 was -1102.3333

$ws.Cells.Item(4,8).Value = 224  # ALC!H4 was 232.5
$ws.Cells.Item(4,9).Value = 248.8  # ALC!I4 was 259
$ws.Cells.Item(4,11).Value = 248.8  # ALC!K4 was 259
$ws.Cells.Item(4,13).Value = -134.8  # ALC!M4 was -145

$ws.Cells.Item(15,8).Value = 349.98  # ALC!H15 was 365.48
$ws.Cells.Item(15,9).Value = 349.98  # ALC!I15 was 365.48
$ws.Cells.Item(15,11).Value = 1049.94  # ALC!K15 was 1096.44
$ws.Cells.Item(15,13).Value = -880.9400000000001  # ALC!M15 was -927.4400000000001

$ws.Cells.Item(41,8).Value = 757.3333  # ALC!H41 was 939
$ws.Cells.Item(41,9).Value = 695.6667  # ALC!I41 was 945.6667
$ws.Cells.Item(41,10).Value = 880.6667  # ALC!J41 was 899
$ws.Cells.Item(41,11).Value = 695.6667  # ALC!K41 was 945.6667
$ws.Cells.Item(41,12).Value = 880.6667  # ALC!L41 was 899
$ws.Cells.Item(41,13).Value = -255.6667  # ALC!M41 was -505.6667
$ws.Cells.Item(41,14).Value = -1760.6667  # ALC!N41 was -1779

$ws.Cells.Item(116,8).Value = 10341.667  # ALC!H116 was 0
$ws.Cells.Item(116,9).Value = 3000  # ALC!I116 was 0
$ws.Cells.Item(116,10).Value = 14012.5  # ALC!J116 was 0
$ws.Cells.Item(116,11).Value = 3000  # ALC!K116 was 0
$ws.Cells.Item(116,12).Value = 14012.5  # ALC!L116 was 0
$ws.Cells.Item(116,13).Value = 442  # ALC!M116 was None
$ws.Cells.Item(116,14).Value = -20896.5  # ALC!N116 was None

$ws.Cells.Item(127,8).Value = 8652.571  # ALC!H127 was 8595.375
$ws.Cells.Item(127,9).Value = 8652.571  # ALC!I127 was 8595.375
$ws.Cells.Item(127,11).Value = 25957.713  # ALC!K127 was 25786.125
$ws.Cells.Item(127,13).Value = -20997.713  # ALC!M127 was -20826.125

$ws.Cells.Item(141,8).Value = 4452  # ALC!H141 was 4063.6667
$ws.Cells.Item(141,9).Value = 1647.1333  # ALC!I141 was 1497.4117
$ws.Cells.Item(141,11).Value = 4941.3999  # ALC!K141 was 4492.2351
$ws.Cells.Item(141,13).Value = 238.6000999999997  # ALC!M141 was 687.7649000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74,8).Value = 561.871  # ARM!H74 was 575.3
$ws.Cells.Item(74,9).Value = 408.10715  # ARM!I74 was 417.33334
$ws.Cells.Item(74,11).Value = 408.10715  # ARM!K74 was 417.33334
$ws.Cells.Item(74,13).Value = 465.89285  # ARM!M74 was 456.66666

$ws.Cells.Item(77,8).Value = 561.871  # ARM!H77 was 575.3
$ws.Cells.Item(77,9).Value = 408.10715  # ARM!I77 was 417.33334
$ws.Cells.Item(77,11).Value = 2040.53575  # ARM!K77 was 2086.6667
$ws.Cells.Item(77,13).Value = 2327.46425  # ARM!M77 was 2281.3333

$ws.Cells.Item(97,8).Value = 589.9231  # ARM!H97 was 486.125
$ws.Cells.Item(97,9).Value = 564.25  # ARM!I97 was 458.66666
$ws.Cells.Item(97,11).Value = 564.25  # ARM!K97 was 458.66666
$ws.Cells.Item(97,13).Value = -68.25  # ARM!M97 was 37.33334000000002

$ws.Cells.Item(102,8).Value = 7938907  # ARM!H102 was 8549471
$ws.Cells.Item(102,9).Value = 13890687  # ARM!I102 was 15874846
$ws.Cells.Item(102,11).Value = 13890687  # ARM!K102 was 15874846
$ws.Cells.Item(102,13).Value = -13889065  # ARM!M102 was -15873224

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80,8).Value = 289.55554  # BSM!H80 was 280.4
$ws.Cells.Item(80,9).Value = 261.75  # BSM!I80 was 239.4
$ws.Cells.Item(80,10).Value = 311.8  # BSM!J80 was 321.4
$ws.Cells.Item(80,11).Value = 261.75  # BSM!K80 was 239.4
$ws.Cells.Item(80,12).Value = 311.8  # BSM!L80 was 321.4
$ws.Cells.Item(80,13).Value = 736.25  # BSM!M80 was 758.6
$ws.Cells.Item(80,14).Value = -2307.8  # BSM!N80 was -2317.4

$ws.Cells.Item(83,8).Value = 289.55554  # BSM!H83 was 280.4
$ws.Cells.Item(83,9).Value = 261.75  # BSM!I83 was 239.4
$ws.Cells.Item(83,10).Value = 311.8  # BSM!J83 was 321.4
$ws.Cells.Item(83,11).Value = 1308.75  # BSM!K83 was 1197
$ws.Cells.Item(83,12).Value = 1559  # BSM!L83 was 1607
$ws.Cells.Item(83,13).Value = 3683.25  # BSM!M83 was 3795
$ws.Cells.Item(83,14).Value = -11543  # BSM!N83 was -11591

$ws.Cells.Item(94,8).Value = 1150.5483  # BSM!H94 was 1090.6666
$ws.Cells.Item(94,9).Value = 671.0833  # BSM!I94 was 631.96155
$ws.Cells.Item(94,11).Value = 671.0833  # BSM!K94 was 631.96155
$ws.Cells.Item(94,13).Value = -220.0833  # BSM!M94 was -180.96155

$ws.Cells.Item(99,8).Value = 1127.2  # BSM!H99 was 1148.75
$ws.Cells.Item(99,9).Value = 1007.86957  # BSM!I99 was 1025.9546
$ws.Cells.Item(99,11).Value = 1007.86957  # BSM!K99 was 1025.9546
$ws.Cells.Item(99,13).Value = 490.13043  # BSM!M99 was 472.0454

$ws.Cells.Item(107,8).Value = 1392.0869  # BSM!H107 was 1342.2084
$ws.Cells.Item(107,10).Value = 1810.5  # BSM!J107 was 1631
$ws.Cells.Item(107,12).Value = 1810.5  # BSM!L107 was 1631
$ws.Cells.Item(107,14).Value = -5650.5  # BSM!N107 was -5471

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31,8).Value = 5019.9414  # CRP!H31 was 5249.933
$ws.Cells.Item(31,9).Value = 3981.375  # CRP!I31 was 4210.1665
$ws.Cells.Item(31,11).Value = 3981.375  # CRP!K31 was 4210.1665
$ws.Cells.Item(31,13).Value = -3686.375  # CRP!M31 was -3915.1665

$ws.Cells.Item(34,8).Value = 5019.9414  # CRP!H34 was 5249.933
$ws.Cells.Item(34,9).Value = 3981.375  # CRP!I34 was 4210.1665
$ws.Cells.Item(34,11).Value = 3981.375  # CRP!K34 was 4210.1665
$ws.Cells.Item(34,13).Value = -3779.375  # CRP!M34 was -4008.1665

$ws.Cells.Item(60,8).Value = 14061.23  # CRP!H60 was 14292.462
$ws.Cells.Item(60,10).Value = 22219  # CRP!J60 was 22970.5
$ws.Cells.Item(60,12).Value = 22219  # CRP!L60 was 22970.5
$ws.Cells.Item(60,14).Value = -23241  # CRP!N60 was -23992.5

$ws.Cells.Item(107,8).Value = 11906469  # CRP!H107 was 31250380
$ws.Cells.Item(107,9).Value = 16668494  # CRP!I107 was 35714692
$ws.Cells.Item(107,10).Value = 1408.9166  # CRP!J107 was 204.5
$ws.Cells.Item(107,11).Value = 16668494  # CRP!K107 was 35714692
$ws.Cells.Item(107,12).Value = 1408.9166  # CRP!L107 was 204.5
$ws.Cells.Item(107,13).Value = -16666574  # CRP!M107 was -35712772
$ws.Cells.Item(107,14).Value = -5248.9166  # CRP!N107 was -4044.5

$ws.Cells.Item(132,8).Value = 1393  # CRP!H132 was 1135
$ws.Cells.Item(132,9).Value = 1393  # CRP!I132 was 1135
$ws.Cells.Item(132,11).Value = 4179  # CRP!K132 was 3405
$ws.Cells.Item(132,13).Value = -1649  # CRP!M132 was -875

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(17,8).Value = 1841.6666  # CUL!H17 was 1883.3334
$ws.Cells.Item(17,9).Value = 200  # CUL!I17 was 0
$ws.Cells.Item(17,10).Value = 2662.5  # CUL!J17 was 1883.3334
$ws.Cells.Item(17,11).Value = 600  # CUL!K17 was 0
$ws.Cells.Item(17,12).Value = 7987.5  # CUL!L17 was 5650.0002
$ws.Cells.Item(17,13).Value = -431  # CUL!M17 was None
$ws.Cells.Item(17,14).Value = -8325.5  # CUL!N17 was -5988.0002

$ws.Cells.Item(39,8).Value = 667  # CUL!H39 was 1055.5
$ws.Cells.Item(39,10).Value = 2996  # CUL!J39 was 2997
$ws.Cells.Item(39,12).Value = 8988  # CUL!L39 was 8991
$ws.Cells.Item(39,14).Value = -9576  # CUL!N39 was -9579

$ws.Cells.Item(75,8).Value = 766.3333  # CUL!H75 was 542.5
$ws.Cells.Item(75,9).Value = 0  # CUL!I75 was 185
$ws.Cells.Item(75,10).Value = 766.3333  # CUL!J75 was 900
$ws.Cells.Item(75,11).Value = 0  # CUL!K75 was 555
$ws.Cells.Item(75,12).ClearContents()  # CUL!L75 remove (was 2700)
$ws.Cells.Item(75,13).Value = 2298.9999  # CUL!M75 was 443
$ws.Cells.Item(75,14).Value = -4294.9999  # CUL!N75 was -4696

$ws.Cells.Item(78,8).Value = 766.3333  # CUL!H78 was 542.5
$ws.Cells.Item(78,9).Value = 0  # CUL!I78 was 185
$ws.Cells.Item(78,10).Value = 766.3333  # CUL!J78 was 900
$ws.Cells.Item(78,11).Value = 0  # CUL!K78 was 1665
$ws.Cells.Item(78,12).ClearContents()  # CUL!L78 remove (was 8100)
$ws.Cells.Item(78,13).Value = 6896.9997  # CUL!M78 was 3327
$ws.Cells.Item(78,14).Value = -16880.9997  # CUL!N78 was -18084

$ws.Cells.Item(92,8).Value = 230.77777  # CUL!H92 was 234.27272
$ws.Cells.Item(92,9).Value = 217.2  # CUL!I92 was 226.57143
$ws.Cells.Item(92,11).Value = 651.5999999999999  # CUL!K92 was 679.71429
$ws.Cells.Item(92,13).Value = 596.4000000000001  # CUL!M92 was 568.28571

$ws.Cells.Item(103,8).Value = 70.25  # CUL!H103 was 106.2
$ws.Cells.Item(103,10).Value = 0  # CUL!J103 was 250
$ws.Cells.Item(103,12).Value = 0  # CUL!L103 was 750
$ws.Cells.Item(103,14).ClearContents()  # CUL!N103 remove (was -2508)

$ws.Cells.Item(120,8).Value = 9095.454  # CUL!H120 was 9605
$ws.Cells.Item(120,10).Value = 10777.777  # CUL!J120 was 11625
$ws.Cells.Item(120,12).Value = 32333.331  # CUL!L120 was 34875
$ws.Cells.Item(120,14).Value = -42009.331  # CUL!N120 was -44551

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(122,8).Value = 40623.27  # GSM!H122 was 39139.258
$ws.Cells.Item(122,9).Value = 1873.2273  # GSM!I122 was 1815.9131
$ws.Cells.Item(122,11).Value = 5619.6819  # GSM!K122 was 5447.7393
$ws.Cells.Item(122,13).Value = -3169.6819  # GSM!M122 was -2997.7393

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40,8).Value = 2017.625  # LTW!H40 was 2145.2856
$ws.Cells.Item(40,9).Value = 2005.8572  # LTW!I40 was 2152.8333
$ws.Cells.Item(40,11).Value = 2005.8572  # LTW!K40 was 2152.8333
$ws.Cells.Item(40,13).Value = -1869.8572  # LTW!M40 was -2016.8333

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(3,8).Value = 11155.25  # WVR!H3 was 14602.333
$ws.Cells.Item(3,9).Value = 1803  # WVR!I3 was 1901.5
$ws.Cells.Item(3,10).Value = 20507.5  # WVR!J3 was 40004
$ws.Cells.Item(3,11).Value = 1803  # WVR!K3 was 1901.5
$ws.Cells.Item(3,12).Value = 20507.5  # WVR!L3 was 40004
$ws.Cells.Item(3,13).Value = -1689  # WVR!M3 was -1787.5
$ws.Cells.Item(3,14).Value = -20735.5  # WVR!N3 was -40232

$ws.Cells.Item(23,8).Value = 550.6  # WVR!H23 was 594.5
$ws.Cells.Item(23,9).Value = 326.5  # WVR!I23 was 549.44446
$ws.Cells.Item(23,10).Value = 700  # WVR!J23 was 1000
$ws.Cells.Item(23,11).Value = 326.5  # WVR!K23 was 549.44446
$ws.Cells.Item(23,12).Value = 700  # WVR!L23 was 1000
$ws.Cells.Item(23,13).Value = -97.5  # WVR!M23 was -320.44446
$ws.Cells.Item(23,14).Value = -1158  # WVR!N23 was -1458

$ws.Cells.Item(81,8).Value = 10661.111  # WVR!H81 was 6364.615
$ws.Cells.Item(81,9).Value = 2725.5  # WVR!I81 was 967.5
$ws.Cells.Item(81,10).Value = 12928.429  # WVR!J81 was 15000
$ws.Cells.Item(81,11).Value = 5451  # WVR!K81 was 1935
$ws.Cells.Item(81,12).Value = 25856.858  # WVR!L81 was 30000
$ws.Cells.Item(81,13).Value = -4390  # WVR!M81 was -874
$ws.Cells.Item(81,14).Value = -27978.858  # WVR!N81 was -32122

$ws.Cells.Item(84,8).Value = 10661.111  # WVR!H84 was 6364.615
$ws.Cells.Item(84,9).Value = 2725.5  # WVR!I84 was 967.5
$ws.Cells.Item(84,10).Value = 12928.429  # WVR!J84 was 15000
$ws.Cells.Item(84,11).Value = 27255  # WVR!K84 was 9675
$ws.Cells.Item(84,12).Value = 129284.29  # WVR!L84 was 150000
$ws.Cells.Item(84,13).Value = -21951  # WVR!M84 was -4371

